# The deck currently has its applied design ("Integral" / Red Violet
# colour scheme) wired to ppt/theme/theme2.xml (used by the slide master
# and by the presentation itself), while ppt/theme/theme1.xml still holds
# the original default "Office Theme" colour scheme (only referenced by
# the notes master). The authored change swaps the two palettes so the
# presentation's applied theme becomes the plain "Office" colours.
#
# Drive this the same way a user would from Design > Variants > Colors >
# Customize Colors: rewrite each of the twelve theme colour slots
# (Background/Text 1-2, Accent 1-6, Hyperlink, Followed Hyperlink) on the
# presentation's theme colour scheme to the stock Office palette.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$themeColors = $master.Theme.ThemeColorScheme

# Order matches ThemeColorScheme.Item(1..12): dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink. Values are OLE RGB (0xBBGGRR) decimal
# equivalents of the stock Office theme palette.
$officeRgb = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officeRgb[$i - 1]
}
